$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update review_text (column F) values for rows 2, 3, 4, 6 with cleaned-up text
$ws.Range("F2").Value = "menarik. tebal. baik. tempat minyak nya bagus bahannya tebal kualitas nya baik muat 2 liter minyak"
$ws.Range("F3").Value = "pas. bagus. kuat dan kokoh"
$ws.Range("F4").Value = "pas sesuai gambar. ok tebal. produk sangat kuat"
$ws.Range("F6").Value = "menarik dan modern. menarik dan ceria. lumayan tebal. Alhamdulillah barang nya sdh sampai dengan baik dan benar,, real picture bagus mangkuk sambel nya,,untuk harga segitu mah wort it,,terima kasih"

# Remove the last product row (row 7) entirely - it is no longer part of the scraped data
$ws.Rows.Item(7).Delete()
